# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 216278f8-bf8a-4779-87dc-fa8dfc2c6d5f entry (row 3) on both
# the zh-cn and de-de report sheets, as part of regenerating the handback
# status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-13 02:40:30"
$wsZhCn.Range("H3").Value = "2016-03-13 02:40:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-13 02:40:34"
$wsDeDe.Range("H3").Value = "2016-03-13 02:40:53"
